# "Printing working as intended now"
# Append the next block of log entries (8/30, OSG building AV Shutdown /
# Demo / Pickup Mic tasks) to the "Logs" sheet, starting at row 328
# (leaving the existing 325-327 gap as-is, matching the sheet's usual
# spacing convention), and add the two brand-new "Special
# Instructions/Comments" / "Time" strings that come with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Each tuple is: row, Staff/TaskType(A), Date(B), Time(C), Building(D), Room(E), Comments(F)
$rows = @(
    @(328, "AV Shutdown", 42612, "1630", "YL",  "280N", $null),
    @(329, "Pickup Mic",  42612, "1630", "YL",  "280N", "Return mic (IR) to KT 516 and place battery in charger"),
    @(330, "AV Shutdown", 42612, "1630", "OSG", "2008", $null),
    @(331, "AV Shutdown", 42612, "1630", "OSG", "1001", $null),
    @(332, "AV Shutdown", 42612, "1630", "OSG", "1002", $null),
    @(333, "AV Shutdown", 42612, "1730", "OSG", "2001", $null),
    @(334, "AV Shutdown", 42612, "1630", "OSG", "2002", $null),
    @(335, "AV Shutdown", 42612, "1630", "OSG", "1005", $null),
    @(336, "Pickup Mic",  42612, "1630", "OSG", "1005", "Return podium mic and 2 desk mics w cables and stands to booth behind stage"),
    @(337, "Demo",        42612, "1630", "OSG", "2004", $null),
    @(338, "Demo",        42612, "1630", "OSG", "2009", $null),
    @(339, "Demo",        42612, "1630", "OSG", "2028", $null),
    @(340, "Demo",        42612, "1630", "OSG", "1003", $null),
    @(341, "AV Shutdown", 42612, "1830", "OSG", "1003", $null),
    @(342, "AV Shutdown", 42612, "1830", "OSG", "2004", $null),
    @(343, "AV Shutdown", 42612, "1830", "OSG", "2009", $null),
    @(344, "AV Shutdown", 42612, "1830", "OSG", "2028", $null),
    @(345, "Demo",        42612, "1830", "OSG", "1004", $null),
    @(346, "Demo",        42612, "1900", "OSG", "1002", $null),
    @(347, "Demo",        42612, "1900", "OSG", "1008", $null),
    @(348, "Demo",        42612, "1900", "OSG", "2001", $null),
    @(349, "AV Shutdown", 42612, "2030", "OSG", "1004", $null),
    @(350, "AV Shutdown", 42612, "2100", "OSG", "1002", $null),
    @(351, "AV Shutdown", 42612, "2100", "OSG", "1008", $null),
    @(352, "AV Shutdown", 42612, "2200", "OSG", "2001", $null)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    if ($r[6] -ne $null) {
        $ws.Range("F$rowNum").Value = $r[6]
    }
}

# Row 336 wraps to two lines in the real workbook (long comment in F336).
$ws.Rows.Item(336).RowHeight = 30

# Leave the view parked where the author left it: scrolled down so the
# new block is visible, with the last-touched cell selected.
$ws.Range("E352").Select()
